$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark (it currently sits in the empty
#    paragraph under the "Arafat Hossain" e-mail line on the title page).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2. Give the paragraph mark of the final "Design and Implementation
#    constraints" paragraph the same Arial / black-color run formatting that
#    its runs already carry. We rebuild the paragraph's own OOXML (runs kept
#    byte for byte) with the two extra <w:rPr> children added to its <w:pPr>,
#    then push it back in with InsertXML so only that paragraph is touched.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastRange = $lastPara.Range.Duplicate

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 wp14">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$constraintParaXml = '<w:body><w:p w:rsidR="00FC71FC" w:rsidRPr="00C93D27" w:rsidRDefault="00FC71FC" w:rsidP="00864E0D"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>The constraint with the implementation is that when the customers create an appointment</w:t></w:r><w:r w:rsidR="007C1EBD"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> how will the system manage these appointments for the doctors because they may have prior appointments from other </w:t></w:r><w:r w:rsidR="00781871"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>sources</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">, one solution is to create a system for the doctors to manage all their appointments from all sources, or just leave the management portion to them.   </w:t></w:r></w:p></w:body>'

$lastRange.InsertXML($pkgHeader + $constraintParaXml + $pkgFooter)

# ---------------------------------------------------------------------------
# 3. Append two brand-new, run-less paragraphs at the very end of the
#    document (the new "Cover Page" spacer + the relocated _GoBack marker).
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)

$newParasXml = '<w:body>' +
  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
  '</w:body>'

$endRange.InsertXML($pkgHeader + $newParasXml + $pkgFooter)

Write-Output "applied"
